# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 201 (pushing the existing rows
# 201..253 down to 202..254) on the single data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 201; Excel shifts rows 201-253 down
# to 202-254 and extends the used range/dimension to A1:T254 automatically.
$ws.Rows(201).Insert()

# Populate the newly inserted row 201 with the new weekly record.
$ws.Cells.Item(201, 1).Value = 6
$ws.Cells.Item(201, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(201, 3).Value = "Metropolitana"
$ws.Cells.Item(201, 4).Value = 44932
$ws.Cells.Item(201, 5).Value = 13
$ws.Cells.Item(201, 6).Value = "Fruta"
$ws.Cells.Item(201, 7).Value = 100101
$ws.Cells.Item(201, 8).Value = "Berries"
$ws.Cells.Item(201, 9).Value = 100101004
$ws.Cells.Item(201, 10).Value = "Frambuesa"
$ws.Cells.Item(201, 11).Value = "Sin especificar"
$ws.Cells.Item(201, 12).Value = "Especial"
$ws.Cells.Item(201, 13).Value = 250
$ws.Cells.Item(201, 14).Value = 8000
$ws.Cells.Item(201, 15).Value = 8000
$ws.Cells.Item(201, 16).Value = 8000
$ws.Cells.Item(201, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(201, 18).Value = "Región del Maule"
$ws.Cells.Item(201, 19).Value = 4000
$ws.Cells.Item(201, 20).Value = 2
